# Update odds figures on Sheet1 (Jogos da Semana FlashScore 2025-01-30)
# per the latest FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2.75
$ws.Range("I2").Value = 2.75
$ws.Range("K2").Value = 1.95
$ws.Range("G3").Value = 1.57
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("S3").Value = 1.98
$ws.Range("T3").Value = 1.88
$ws.Range("W3").Value = 3.4
$ws.Range("X3").Value = 1.3
$ws.Range("G4").Value = 2.05
$ws.Range("J4").Value = 2.88
$ws.Range("K4").Value = 1.95
$ws.Range("T4").Value = 1.5
$ws.Range("G5").Value = 2.4
$ws.Range("K5").Value = 1.91
$ws.Range("T5").Value = 1.44
$ws.Range("H9").Value = 3.05
$ws.Range("I9").Value = 2.3
$ws.Range("J9").Value = 3.5
$ws.Range("K9").Value = 2.07
$ws.Range("L9").Value = 2.82
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 3
$ws.Range("S9").Value = 1.87
$ws.Range("T9").Value = 1.75
$ws.Range("W9").Value = 2.95
$ws.Range("X9").Value = 1.3
$ws.Range("AA9").Value = 1.65
$ws.Range("AB9").Value = 1.98
$ws.Range("AC9").Value = 9.5
$ws.Range("AD9").Value = 16.5
$ws.Range("AG9").Value = 26
$ws.Range("AH9").Value = 32
$ws.Range("AI9").Value = 9.25
$ws.Range("AJ9").Value = 6
$ws.Range("AO9").Value = 11.5
$ws.Range("AQ9").Value = 24
$ws.Range("AR9").Value = 18.5
$ws.Range("AS9").Value = 27
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8
$ws.Range("O10").Value = 1.4
$ws.Range("S10").Value = 2.25
$ws.Range("T10").Value = 1.62
$ws.Range("W10").Value = 4.33
$ws.Range("X10").Value = 1.2
$ws.Range("G11").Value = 2.15
$ws.Range("H11").Value = 2.9
$ws.Range("I11").Value = 3.3
$ws.Range("J11").Value = 3.1
$ws.Range("K11").Value = 1.91
$ws.Range("L11").Value = 4.33
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 6.5
$ws.Range("O11").Value = 1.5
$ws.Range("Q11").Value = 1.93
$ws.Range("R11").Value = 1.88
$ws.Range("X11").Value = 1.17
$ws.Range("Y11").Value = 1.57
$ws.Range("Z11").Value = 2.25
$ws.Range("AA11").Value = 2.2
$ws.Range("AB11").Value = 1.62
$ws.Range("AD11").Value = 9.5
$ws.Range("AE11").Value = 10
$ws.Range("AF11").Value = 21
$ws.Range("AI11").Value = 6.5
$ws.Range("AN11").Value = 8
$ws.Range("AO11").Value = 15
$ws.Range("G13").Value = 5
$ws.Range("H13").Value = 4.55
$ws.Range("J13").Value = 4.55
$ws.Range("K13").Value = 2.72
$ws.Range("L13").Value = 1.91
$ws.Range("P13").Value = 6.3
$ws.Range("S13").Value = 1.29
$ws.Range("X13").Value = 2.02
$ws.Range("Y13").Value = 1.17
$ws.Range("Z13").Value = 4.35
$ws.Range("AA13").Value = 1.38
$ws.Range("AB13").Value = 2.8
$ws.Range("AC13").Value = 32
$ws.Range("AD13").Value = 50
$ws.Range("AE13").Value = 17.5
$ws.Range("AF13").Value = 100
$ws.Range("AL13").Value = 29
$ws.Range("AQ13").Value = 14
$ws.Range("AR13").Value = 10.5
$ws.Range("G15").Value = 1.9
$ws.Range("H15").Value = 3.5
$ws.Range("I15").Value = 3.8
$ws.Range("T15").Value = 2
$ws.Range("AD15").Value = 9.5
$ws.Range("AJ15").Value = 7
$ws.Range("AK15").Value = 15
$ws.Range("AO15").Value = 21
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 8.5
$ws.Range("I16").Value = 1.18
$ws.Range("J16").Value = 9
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1.5
$ws.Range("M16").Value = 1.01
$ws.Range("N16").Value = 29
$ws.Range("O16").Value = 1.08
$ws.Range("P16").Value = 7
$ws.Range("S16").Value = 1.3
$ws.Range("T16").Value = 3.4
$ws.Range("W16").Value = 1.8
$ws.Range("X16").Value = 1.91
$ws.Range("Y16").Value = 1.2
$ws.Range("Z16").Value = 4.33
$ws.Range("AC16").Value = 34
$ws.Range("AE16").Value = 29
$ws.Range("AF16").Value = 126
$ws.Range("AI16").Value = 29
$ws.Range("AJ16").Value = 17
$ws.Range("AM16").Value = 600
$ws.Range("AN16").Value = 12
$ws.Range("AO16").Value = 8
$ws.Range("AP16").Value = 11
$ws.Range("G19").Value = 1.55
$ws.Range("H19").Value = 4.15
$ws.Range("I19").Value = 5.1
$ws.Range("J19").Value = 2.05
$ws.Range("K19").Value = 2.42
$ws.Range("L19").Value = 4.9
$ws.Range("AA19").Value = 1.6
$ws.Range("AB19").Value = 2.18
$ws.Range("AD19").Value = 8.75
$ws.Range("AF19").Value = 12
$ws.Range("AG19").Value = 11.25
$ws.Range("AH19").Value = 19.5
$ws.Range("AK19").Value = 14
$ws.Range("AL19").Value = 50
$ws.Range("AM19").Value = 300
$ws.Range("AN19").Value = 19
$ws.Range("AP19").Value = 16
$ws.Range("AQ19").Value = 90
$ws.Range("AS19").Value = 37
